$d = $word.ActiveDocument

# Change 1: insert clarification about the exponential distribution of the
# deterioration time, right after "un processo di deterioramento" and before
# the following period. Replace the whole sentence (spanning the original
# three runs) so the paragraph collapses to a single run, matching the diff.
$d.Content.Find.Execute(
    "Il blocco b. della banca degli organi disponibili per l" + [char]0x2019 + "assegnamento è modellato come 4 code FIFO, una per ogni gruppo sanguigno. Gli organi una volta entrati nel sistema sono soggetti ad un processo di deterioramento. Pertanto una volta entrati nel sistema gli organi possono essere associati ad un paziente compatibile e quindi essere trapiantati, oppure possono deteriorarsi a tal punto da non poter essere utilizzati per un trapianto in quanto potrebbero provocare reazioni nel paziente recipienti.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Il blocco b. della banca degli organi disponibili per l" + [char]0x2019 + "assegnamento è modellato come 4 code FIFO, una per ogni gruppo sanguigno. Gli organi una volta entrati nel sistema sono soggetti ad un processo di deterioramento (con distribuzione esponenziale del tempo di deterioramento). Pertanto una volta entrati nel sistema gli organi possono essere associati ad un paziente compatibile e quindi essere trapiantati, oppure possono deteriorarsi a tal punto da non poter essere utilizzati per un trapianto in quanto potrebbero provocare reazioni nel paziente recipienti.",
    2
)

# Change 2: no wording change, but the three runs of this bullet point get
# merged into a single run (same formatting throughout), matching the diff.
$d.Content.Find.Execute(
    "Nel modello non sono stati considerati tipi diversi di organi, ma un unico tipo di organo generico.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nel modello non sono stati considerati tipi diversi di organi, ma un unico tipo di organo generico.",
    2
)
